# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.723.35'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '2.329.89'
$ws.Range('E3').Value = '  +4.55%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''271.16'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').Value = '''95.62'
$ws.Range('E6').Value = '  +2.48%  '
$ws.Range('E7').Value = '  +0.83%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = '''0.623'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('D10').Value = '''45.42'
$ws.Range('E10').Value = '  -2.41%  '
$ws.Range('D11').Value = '''0.0945'
$ws.Range('E11').Value = '  +2.98%  '
$ws.Range('D12').Value = '''8.11'
$ws.Range('E12').Value = '  +1.08%  '
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').Value = '2.675.61'
$ws.Range('E14').Value = '  +4.30%  '
$ws.Range('D15').Value = '''15.65'
$ws.Range('E15').Value = '  +3.80%  '
$ws.Range('D16').Value = '''0.864'
$ws.Range('E16').Value = '  +8.28%  '
$ws.Range('D17').Value = '2.342.05'
$ws.Range('E17').Value = '  +4.44%  '
$ws.Range('D18').Value = '43.710.75'
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range('E19').Value = '  +5.97%  '
$ws.Range('D20').Value = '''6.42'
$ws.Range('E20').Value = '  +7.17%  '
$ws.Range('D21').Value = '''72.63'
$ws.Range('E21').Value = '  +3.28%  '
$ws.Range('D22').Value = '''239.36'
$ws.Range('E22').Value = '  +3.10%  '
$ws.Range('D23').Value = '''2.28'
$ws.Range('E23').Value = '  -2.17%  '
$ws.Range('D24').Value = '''9.36'
$ws.Range('E24').Value = '  +7.31%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  +1.82%  '
$ws.Range('E27').Value = '  +1.74%  '
$ws.Range('D28').Value = '''3.48'
$ws.Range('E28').Value = '  -2.10%  '
$ws.Range('D29').Value = '''2.27'
$ws.Range('E29').Value = '  +0.38%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '''38.30'
$ws.Range('E30').Value = '  -2.90%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '''22.56'
$ws.Range('E31').Value = '  +8.77%  '
$ws.Range('D32').Value = '''172.99'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').Value = '''0.0901'
$ws.Range('E33').Value = '  -2.54%  '
$ws.Range('D34').Value = '''5.50'
$ws.Range('E34').Value = '  +1.42%  '
$ws.Range('E35').Value = '  +2.62%  '
$ws.Range('E36').Value = '  +4.00%  '
$ws.Range('E37').Value = '  -2.37%  '
$ws.Range('D38').Value = '''4.38'
$ws.Range('E38').Value = '  +2.30%  '
$ws.Range('E39').Value = '  -3.33%  '
$ws.Range('E40').Value = '  +10.25%  '
$ws.Range('E41').Value = '  +9.82%  '
$ws.Range('D42').Value = '''1.38'
$ws.Range('E42').Value = '  +20.73%  '
$ws.Range('D43').Value = '''12.15'
$ws.Range('E43').Value = '  -3.05%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''9.22'
$ws.Range('E44').Value = '  +10.00%  '
$ws.Range('B45').Value = 'MultiversX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D45').Value = '''62.34'
$ws.Range('E45').Value = '  -0.66%  '
$ws.Range('D46').Value = '''5.36'
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('E47').Value = '  +4.96%  '
$ws.Range('D48').Value = '''100.48'
$ws.Range('E48').Value = '  +0.77%  '
$ws.Range('E49').Value = '  +1.84%  '
$ws.Range('D50').Value = '''0.189'
$ws.Range('E50').Value = '  +16.71%  '
$ws.Range('D51').Value = '2.553.57'
$ws.Range('E51').Value = '  +4.07%  '
